$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 488.75
$ws.Range("I33").Value = 447.85715
$ws.Range("J33").Value = 775
$ws.Range("K33").Value = 447.85715
$ws.Range("L33").Value = 775
$ws.Range("M33").Value = -218.85715
$ws.Range("N33").Value = -1233

$ws.Range("H70").Value = 2616.6667
$ws.Range("I70").Value = 2500
$ws.Range("J70").Value = 2640
$ws.Range("K70").Value = 7500
$ws.Range("L70").Value = 7920
$ws.Range("M70").Value = -7230
$ws.Range("N70").Value = -8460

$ws.Range("H73").Value = 2616.6667
$ws.Range("I73").Value = 2500
$ws.Range("J73").Value = 2640
$ws.Range("K73").Value = 7500
$ws.Range("L73").Value = 7920
$ws.Range("M73").Value = -6564
$ws.Range("N73").Value = -9792

$ws.Range("H98").Value = 851.7895
$ws.Range("I98").Value = 924.25
$ws.Range("K98").Value = 924.25
$ws.Range("M98").Value = 573.75

$ws.Range("H122").Value = 851.7895
$ws.Range("I122").Value = 924.25
$ws.Range("K122").Value = 2772.75
$ws.Range("M122").Value = -322.75

$ws.Range("H137").Value = 2080.5908
$ws.Range("I137").Value = 1592.5883
$ws.Range("K137").Value = 4777.7649
$ws.Range("M137").Value = -2227.7649

$ws.Range("H138").Value = 3291
$ws.Range("J138").Value = 4274.75
$ws.Range("L138").Value = 12824.25
$ws.Range("N138").Value = -23104.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 196.91667
$ws.Range("I5").Value = 209.28572
$ws.Range("J5").Value = 179.6
$ws.Range("K5").Value = 209.28572
$ws.Range("L5").Value = 179.6
$ws.Range("M5").Value = -97.28572
$ws.Range("N5").Value = -403.6

$ws.Range("H8").Value = 4801800
$ws.Range("J8").Value = 3000
$ws.Range("L8").Value = 3000
$ws.Range("N8").Value = -3288

$ws.Range("H11").Value = 1000274.75
$ws.Range("J11").Value = 366.33334
$ws.Range("L11").Value = 366.33334
$ws.Range("N11").Value = -654.33334

$ws.Range("H14").Value = 1280
$ws.Range("J14").Value = 1100
$ws.Range("L14").Value = 1100
$ws.Range("N14").Value = -1450

$ws.Range("H32").Value = 5090.5186
$ws.Range("I32").Value = 3065.682
$ws.Range("K32").Value = 3065.682
$ws.Range("M32").Value = -2778.682

$ws.Range("H117").Value = 35000
$ws.Range("J117").Value = 35000
$ws.Range("L117").Value = 35000
$ws.Range("N117").Value = -44178

$ws.Range("H122").Value = 2552.7
$ws.Range("I122").Value = 2414.875
$ws.Range("J122").Value = 3104
$ws.Range("K122").Value = 7244.625
$ws.Range("L122").Value = 9312
$ws.Range("M122").Value = -4794.625
$ws.Range("N122").Value = -14212

$ws.Range("H132").Value = 974.3333
$ws.Range("I132").Value = 962
$ws.Range("K132").Value = 2886
$ws.Range("M132").Value = -356

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 196.91667
$ws.Range("I4").Value = 209.28572
$ws.Range("J4").Value = 179.6
$ws.Range("K4").Value = 209.28572
$ws.Range("L4").Value = 179.6
$ws.Range("M4").Value = -94.28572
$ws.Range("N4").Value = -409.6

$ws.Range("H12").Value = 1800
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1800
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1800
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -2136

$ws.Range("H94").Value = 4283.4
$ws.Range("I94").Value = 3925.1
$ws.Range("K94").Value = 3925.1
$ws.Range("M94").Value = -3474.1

$ws.Range("H117").Value = 75000
$ws.Range("J117").Value = 75000
$ws.Range("L117").Value = 75000
$ws.Range("N117").Value = -84178

$ws.Range("H134").Value = 3149.3333
$ws.Range("I134").Value = 3379.7
$ws.Range("K134").Value = 10139.1
$ws.Range("M134").Value = -7604.099999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 491.66666
$ws.Range("I16").Value = 537.5
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 537.5
$ws.Range("L16").Value = 400
$ws.Range("M16").Value = -250.5
$ws.Range("N16").Value = -974

$ws.Range("H31").Value = 2829.2222
$ws.Range("J31").Value = 2522.25
$ws.Range("L31").Value = 2522.25
$ws.Range("N31").Value = -3112.25

$ws.Range("H34").Value = 2829.2222
$ws.Range("J34").Value = 2522.25
$ws.Range("L34").Value = 2522.25
$ws.Range("N34").Value = -2926.25

$ws.Range("H113").Value = 491.66666
$ws.Range("I113").Value = 537.5
$ws.Range("J113").Value = 400
$ws.Range("K113").Value = 537.5
$ws.Range("L113").Value = 400
$ws.Range("M113").Value = 1632.5
$ws.Range("N113").Value = -4740

$ws.Range("H122").Value = 840.2857
$ws.Range("I122").Value = 858.7
$ws.Range("K122").Value = 2576.1
$ws.Range("M122").Value = -126.1000000000004

$ws.Range("H134").Value = 1963.3334
$ws.Range("I134").Value = 1963.3334
$ws.Range("K134").Value = 5890.0002
$ws.Range("M134").Value = -3355.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1435.125
$ws.Range("I109").Value = 211.57143
$ws.Range("K109").Value = 634.71429
$ws.Range("M109").Value = 405.28571

$ws.Range("H131").Value = 1084.5385
$ws.Range("J131").Value = 1083.25
$ws.Range("L131").Value = 3249.75
$ws.Range("N131").Value = -13329.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3056
$ws.Range("I80").Value = 2834
$ws.Range("K80").Value = 2834
$ws.Range("M80").Value = -1836

$ws.Range("H83").Value = 3056
$ws.Range("I83").Value = 2834
$ws.Range("K83").Value = 14170
$ws.Range("M83").Value = -9178

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1779.7
$ws.Range("I22").Value = 1166.1666
$ws.Range("J22").Value = 2700
$ws.Range("K22").Value = 1166.1666
$ws.Range("L22").Value = 2700
$ws.Range("M22").Value = -871.1666
$ws.Range("N22").Value = -3290

$ws.Range("H27").Value = 1779.7
$ws.Range("I27").Value = 1166.1666
$ws.Range("J27").Value = 2700
$ws.Range("K27").Value = 1166.1666
$ws.Range("L27").Value = 2700
$ws.Range("M27").Value = -1059.1666
$ws.Range("N27").Value = -2914

$ws.Range("H61").Value = 1780
$ws.Range("I61").Value = 1780
$ws.Range("K61").Value = 1780
$ws.Range("M61").Value = -1578

$ws.Range("H100").Value = 5165.1665
$ws.Range("I100").Value = 4248.5
$ws.Range("J100").Value = 6998.5
$ws.Range("K100").Value = 4248.5
$ws.Range("L100").Value = 6998.5
$ws.Range("M100").Value = -3707.5
$ws.Range("N100").Value = -8080.5

$ws.Range("H113").Value = 1780
$ws.Range("I113").Value = 1780
$ws.Range("K113").Value = 1780
$ws.Range("M113").Value = 390

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 18000
$ws.Range("J22").Value = 18000
$ws.Range("L22").Value = 18000
$ws.Range("N22").Value = -18586

$ws.Range("H100").Value = 4358727
$ws.Range("I100").Value = 7746987.5
$ws.Range("K100").Value = 15493975
$ws.Range("M100").Value = -15493434

$ws.Range("H132").Value = 1779.8
$ws.Range("I132").Value = 974.75
$ws.Range("K132").Value = 2924.25
$ws.Range("M132").Value = -394.25
